# ertac_additional_variables.xlsx
# "Default NH3 added for more states / Uses CT values for defaults"
#
# For every data row (2..158) whose "Generic NH3 Emission Rate lbs/mmbtu"
# (column G) is still the placeholder 0, fill it in with Connecticut's (CT)
# own NH3 rate for the matching ertac_fuel_unit_type_bin (column B), the
# same way CT's real measured values already exist for Boiler Gas, Combined
# Cycle Gas and Oil. Rows whose fuel type has no non-zero CT default (Coal,
# Simple Cycle Gas) - and any row that already has a non-zero value - are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 158
$fuelCol = 2
$nh3Col = 7

# --- Build the CT (Connecticut) fuel-type -> NH3 rate lookup -------------
$ctMap = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $state = $ws.Cells.Item($r, 1).Value()
    if ($state -eq "CT") {
        $fuel = $ws.Cells.Item($r, $fuelCol).Value()
        $rate = $ws.Cells.Item($r, $nh3Col).Value()
        if ($fuel -ne $null) {
            $ctMap[$fuel.ToUpper()] = $rate
        }
    }
}

# --- Apply CT's defaults to every other row still sitting at 0 -----------
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $fuel = $ws.Cells.Item($r, $fuelCol).Value()
    $current = $ws.Cells.Item($r, $nh3Col).Value()
    if ($fuel -ne $null -and $current -eq 0) {
        $key = $fuel.ToUpper()
        if ($ctMap.ContainsKey($key)) {
            $default = $ctMap[$key]
            if ($default -ne 0) {
                $ws.Cells.Item($r, $nh3Col).Value = $default
            }
        }
    }
}

# --- Turn on AutoFilter over the whole table ------------------------------
$ws.Range("A1:P158").AutoFilter()

# AutoFilter in real Excel also creates a hidden workbook-level (sheet
# scoped) defined name that records the filtered range.
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$P`$158")
$filterName.Visible = $false

# --- Restore the active selection to B1 -----------------------------------
$ws.Range("B1").Select()
